$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shrink the table from 10x10 (A1:J10) down to 7x7 (A1:G7) by removing the
# extra columns (H:J) and extra rows (8:10).
$ws.Range("H1:J10").Delete() | Out-Null
$ws.Range("A8:G10").Delete() | Out-Null

# New values for the remaining data matrix (B2:G7)
$v1 = [double]"0.8069333896388119"
$v2 = [double]"0.006070787794139079"
$v3 = [double]"1.708905370183168e-06"

# Row 2 (i1)
$ws.Cells.Item(2,2).Value = 0
$ws.Cells.Item(2,3).Value = 0
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = 0
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(2,7).Value = 0

# Row 3 (i2)
$ws.Cells.Item(3,2).Value = 0
$ws.Cells.Item(3,3).Value = 0
$ws.Cells.Item(3,4).Value = 0
$ws.Cells.Item(3,5).Value = 0
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(3,7).Value = $v1

# Row 4 (i3)
$ws.Cells.Item(4,2).Value = 0
$ws.Cells.Item(4,3).Value = 0
$ws.Cells.Item(4,4).Value = 0
$ws.Cells.Item(4,5).Value = 0
$ws.Cells.Item(4,6).Value = 0
$ws.Cells.Item(4,7).Value = 0

# Row 5 (i4)
$ws.Cells.Item(5,2).Value = 0
$ws.Cells.Item(5,3).Value = 0
$ws.Cells.Item(5,4).Value = 0
$ws.Cells.Item(5,5).Value = 0
$ws.Cells.Item(5,6).Value = 0
$ws.Cells.Item(5,7).Value = $v2

# Row 6 (i5)
$ws.Cells.Item(6,2).Value = 0
$ws.Cells.Item(6,3).Value = 0
$ws.Cells.Item(6,4).Value = 0
$ws.Cells.Item(6,5).Value = 0
$ws.Cells.Item(6,6).Value = 0
$ws.Cells.Item(6,7).Value = $v3

# Row 7 (i6)
$ws.Cells.Item(7,2).Value = 0
$ws.Cells.Item(7,3).Value = $v1
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = $v2
$ws.Cells.Item(7,6).Value = $v3
$ws.Cells.Item(7,7).Value = 0
